$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 3 (trial 2)
$ws.Range("D3").Value = 3
$ws.Range("F3").Value = 3
$ws.Range("H3").Value = 46

# Update the active selection to D3
$ws.Range("D3").Select()
